$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 100.8373843333333
$ws.Range("H2").Value = 302.512153
$ws.Range("I2").Value = 0.6551985585448407
$ws.Range("J2").Value = 0.6551985585448408
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 13.267299
$ws.Range("N2").Value = 39.801897
$ws.Range("O2").Value = 0.4248136128385448
$ws.Range("P2").Value = 0.4248136128385448
$ws.Range("Q2").Value = 1337.839728328249
$ws.Range("R2").Value = 12040.55755495424
$ws.Range("S2").Value = 0.2783372667820406
$ws.Range("T2").Value = 0.2783372667820406

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 100.8373843333333
$ws.Range("H3").Value = 302.512153
$ws.Range("I3").Value = 0.6551985585448407
$ws.Range("J3").Value = 0.6551985585448408
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 10.340832
$ws.Range("N3").Value = 31.022496
$ws.Range("O3").Value = 0.3311093088108164
$ws.Range("P3").Value = 0.3311093088108164
$ws.Range("Q3").Value = 1042.742450710432
$ws.Range("R3").Value = 9384.682056393889
$ws.Range("S3").Value = 0.2169423418536254
$ws.Range("T3").Value = 0.2169423418536255

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 100.8373843333333
$ws.Range("H4").Value = 302.512153
$ws.Range("I4").Value = 0.6551985585448407
$ws.Range("J4").Value = 0.6551985585448408
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 7.622739666666667
$ws.Range("N4").Value = 22.868219
$ws.Range("O4").Value = 0.2440770783506388
$ws.Range("P4").Value = 0.2440770783506388
$ws.Range("Q4").Value = 768.6571294406119
$ws.Range("R4").Value = 6917.914164965508
$ws.Range("S4").Value = 0.1599189499091747
$ws.Range("T4").Value = 0.1599189499091747

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 23.90796933333333
$ws.Range("H5").Value = 71.723908
$ws.Range("I5").Value = 0.1553438454249564
$ws.Range("J5").Value = 0.1553438454249564
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 13.267299
$ws.Range("N5").Value = 39.801897
$ws.Range("O5").Value = 0.4248136128385448
$ws.Range("P5").Value = 0.4248136128385448
$ws.Range("Q5").Value = 317.194177628164
$ws.Range("R5").Value = 2854.747598653476
$ws.Range("S5").Value = 0.06599218020720818
$ws.Range("T5").Value = 0.06599218020720818

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 23.90796933333333
$ws.Range("H6").Value = 71.723908
$ws.Range("I6").Value = 0.1553438454249564
$ws.Range("J6").Value = 0.1553438454249564
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 10.340832
$ws.Range("N6").Value = 31.022496
$ws.Range("O6").Value = 0.3311093088108164
$ws.Range("P6").Value = 0.3311093088108164
$ws.Range("Q6").Value = 247.228294337152
$ws.Range("R6").Value = 2225.054649034368
$ws.Range("S6").Value = 0.05143579328667162
$ws.Range("T6").Value = 0.05143579328667162

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 23.90796933333333
$ws.Range("H7").Value = 71.723908
$ws.Range("I7").Value = 0.1553438454249564
$ws.Range("J7").Value = 0.1553438454249564
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 7.622739666666667
$ws.Range("N7").Value = 22.868219
$ws.Range("O7").Value = 0.2440770783506388
$ws.Range("P7").Value = 0.2440770783506388
$ws.Range("Q7").Value = 182.2442261866502
$ws.Range("R7").Value = 1640.198035679852
$ws.Range("S7").Value = 0.03791587193107661
$ws.Range("T7").Value = 0.03791587193107661

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 29.15819666666667
$ws.Range("H8").Value = 87.47459
$ws.Range("I8").Value = 0.1894575960302029
$ws.Range("J8").Value = 0.1894575960302029
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 13.267299
$ws.Range("N8").Value = 39.801897
$ws.Range("O8").Value = 0.4248136128385448
$ws.Range("P8").Value = 0.4248136128385448
$ws.Range("Q8").Value = 386.85051347747
$ws.Range("R8").Value = 3481.65462129723
$ws.Range("S8").Value = 0.08048416584929605
$ws.Range("T8").Value = 0.08048416584929605

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 29.15819666666667
$ws.Range("H9").Value = 87.47459
$ws.Range("I9").Value = 0.1894575960302029
$ws.Range("J9").Value = 0.1894575960302029
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 10.340832
$ws.Range("N9").Value = 31.022496
$ws.Range("O9").Value = 0.3311093088108164
$ws.Range("P9").Value = 0.3311093088108164
$ws.Range("Q9").Value = 301.52001315296
$ws.Range("R9").Value = 2713.68011837664
$ws.Range("S9").Value = 0.06273117367051936
$ws.Range("T9").Value = 0.06273117367051936

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 29.15819666666667
$ws.Range("H10").Value = 87.47459
$ws.Range("I10").Value = 0.1894575960302029
$ws.Range("J10").Value = 0.1894575960302029
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 7.622739666666667
$ws.Range("N10").Value = 22.868219
$ws.Range("O10").Value = 0.2440770783506388
$ws.Range("P10").Value = 0.2440770783506388
$ws.Range("Q10").Value = 222.2653423394678
$ws.Range("R10").Value = 2000.38808105521
$ws.Range("S10").Value = 0.04624225651038751
$ws.Range("T10").Value = 0.04624225651038751

Write-Output "applied edits"
